# Update "想去人数" (number of people interested) counts for a fresh data pull.
# This mirrors the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 253
$ws.Range("F3").Value  = 2638
$ws.Range("F7").Value  = 1965
$ws.Range("F11").Value = 2436
$ws.Range("F12").Value = 525
$ws.Range("F18").Value = 8974
$ws.Range("F20").Value = 6964
$ws.Range("F21").Value = 11372
$ws.Range("F24").Value = 226
$ws.Range("F25").Value = 314
$ws.Range("F26").Value = 540
$ws.Range("F27").Value = 2496
$ws.Range("F29").Value = 188
$ws.Range("F30").Value = 2408
$ws.Range("F31").Value = 598
$ws.Range("F34").Value = 770
$ws.Range("F35").Value = 332
$ws.Range("F37").Value = 489

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 11

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 135

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 253
$ws.Range("F6").Value  = 2638
$ws.Range("F10").Value = 1965
$ws.Range("F15").Value = 2436
$ws.Range("F17").Value = 525
$ws.Range("F23").Value = 8974
$ws.Range("F25").Value = 6964
$ws.Range("F26").Value = 11372
$ws.Range("F27").Value = 11
$ws.Range("F29").Value = 226
$ws.Range("F30").Value = 314
$ws.Range("F32").Value = 540
$ws.Range("F34").Value = 2496
$ws.Range("F38").Value = 188
$ws.Range("F46").Value = 489

$wb.Save()
